# Livrable du Pré-TPI le 03.04.2020 -- append a new journal entry as row 17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 16) down to the new row 17, so the
# new entry keeps the same borders/number-formats/fill as the rest of the log.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new entry's values.
$ws.Range("A17").Value = 43924
$ws.Range("B17").Value = 1.5
$ws.Range("C17").Value = "Rendu du Pré-TPI"

# Update selection to mirror the authored workbook state
$ws.Range("C17").Select()
